# "codigo do sensor de temperatura" - implementacao da leitura da
# temperatura pelo painel de controle.
#
# Adds a new command row (CMD 127) to the "tabela comandos" sheet that
# documents the new "ler valor da temperatura" command and its "temp
# atual" return value, widens the new RET column (F) to fit the text,
# and leaves the selection where the author left off (C17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New command row describing the temperature-sensor read command.
$ws.Range("C16").Value = 127
$ws.Range("E16").Value = "ler valor da temperatura"
$ws.Range("F16").Value = "temp atual"

# Give column F (the new RET value for this row) a sensible width.
$ws.Columns(6).ColumnWidth = 13.7

# Match the author's final cell selection.
$ws.Range("C17").Select()
